# "fixed a spelling mistake"
#
# Slide 1 (Title Slide): the title "Federal Health Care Programs" was
# previously split across three separate runs ("Federal ", "Health Care ",
# "Programs"). Re-typing the whole title collapses it into a single run.
$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
$titleRange1 = $s1.Shapes.Item(1).TextFrame.TextRange
# Force a real text replacement (round-tripping through a throwaway value
# that shares no prefix/suffix with the target text) so the run gets
# rewritten as a single piece instead of being left untouched because the
# net text happens to already equal the concatenation of the old runs.
$titleRange1.Text = "X"
$titleRange1.Text = "Federal Health Care Programs"

# Slide 2: title said "Introduction to Medicare" but should read
# "Introduction to Medicaid" -- fix the misspelled/wrong word only, which is
# how PowerPoint naturally keeps the untouched leading text as its own run
# and the replaced word as a new run.
$s2 = $p.Slides.Item(2)
$titleRange2 = $s2.Shapes.Item(1).TextFrame.TextRange
$word = $titleRange2.Characters(17, 8)
$word.Text = "Medicaid"
